$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://play.google.com/store/apps/details?id=com.ucmqed.peqqus",
    "https://play.google.com/store/apps/details?id=com.upticker.app.v2",
    "https://play.google.com/store/apps/details?id=com.digitaljungle",
    "https://play.google.com/store/apps/details?id=com.advanced.voiceChanger",
    "https://play.google.com/store/apps/details?id=com.app.malgopay",
    "https://play.google.com/store/apps/details?id=com.machi.call.ChatLiveCall.mochi",
    "https://play.google.com/store/apps/details?id=co.sansa.vhgmw",
    "https://play.google.com/store/apps/details?id=com.astropower.user",
    "https://play.google.com/store/apps/details?id=com.fotoapp.perfectme.bodyshape.photoeditor.bodyeditor",
    "https://play.google.com/store/apps/details?id=com.mb.hexasort.coinstackgame",
    "https://play.google.com/store/apps/details?id=antitheftalarm.donttouchmyphone.antipro",
    "https://play.google.com/store/apps/details?id=dailybend.bend.stretch.flexibility.fitness.exercise.back",
    "https://play.google.com/store/apps/details?id=photo.video.maker.editor",
    "https://play.google.com/store/apps/details?id=pl.idreams.SkyForceReloaded",
    "https://play.google.com/store/apps/details?id=com.fakecall.prankcall.videocall.voicechat",
    "https://play.google.com/store/apps/details?id=com.aamdhane.app",
    "https://play.google.com/store/apps/details?id=com.idramas.tv",
    "https://play.google.com/store/apps/details?id=com.trackmyphone.find.myphone.phonelocator.antitheft",
    "https://play.google.com/store/apps/details?id=com.earthvibe.wall",
    "https://play.google.com/store/apps/details?id=com.eventbeep.android",
    "https://play.google.com/store/apps/details?id=com.wowchat.wallet",
    "https://play.google.com/store/apps/details?id=com.smartearn.gpstracker",
    "https://play.google.com/store/apps/details?id=xyz.multipl.multipl",
    "https://play.google.com/store/apps/details?id=com.live.map.earth.app",
    "https://play.google.com/store/apps/details?id=com.appspelago.peekababy",
    "https://play.google.com/store/apps/details?id=com.qr.cvn.barcode.scanner.fhnn",
    "https://play.google.com/store/apps/details?id=com.happyverse.spellinggame",
    "https://play.google.com/store/apps/details?id=com.prank.crazy.trickly.aos",
    "https://play.google.com/store/apps/details?id=grocery.beldara.fresh",
    "https://play.google.com/store/apps/details?id=com.tnibro",
    "https://play.google.com/store/apps/details?id=com.navigation.offlinemaps.radar.fortravel"
)

$startRow = 81
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $urls[$i]
}
